# "New materials for semester 2 2021-2022"
#
# 1) Refresh the cached "datetimeFigureOut" field text (11/26/2021 -> 2/21/2022)
#    on every slide layout and the slide master's Date placeholder.
# 2) Rewrite slide 1's title + subtitle with the new intro / contact blurb.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholders on the slide master and all slide layouts.
# ---------------------------------------------------------------------------
$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "2/21/2022"
    }
}

for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "2/21/2022"
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 1 title + subtitle text.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)

$slide1.Shapes.Item(1).TextFrame.TextRange.Text = "Giới thiệu"

$subtitle = $slide1.Shapes.Item(2).TextFrame.TextRange
$subtitle.Text = "Vũ Tuấn Hải, Bộ môn Phát triển phần mềm, khoa CNPM`rEmail: haivt@uit.edu.vn`rFacebook: fb.com/vutuanhai237`rTài liệu học tập: moodle hoặc https://www.facebook.com/groups/bht.cnpm.uit"

for ($i = 1; $i -le $subtitle.Paragraphs().Count; $i++) {
    $subtitle.Paragraphs($i).ParagraphFormat.Alignment = 1
}
